$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 9 (Excel copies formatting from the row above by default)
$ws.Rows.Item(9).Insert()

# New row 9 gets the content that used to be on row 8 ("Upload")
$ws.Range("A9").Value = "Upload"
$ws.Range("B9").Value = $false
$ws.Range("C9").Value = $false

# Old row 8 becomes the new "Force" entry
$ws.Range("A8").Value = "Force"

# Update selection / frozen pane to match the recorded post-edit view state
$ws.Range("A9").Select()
